# motilal_portfolio_change_engine: insert the "Industry" column (new col C)
# into the Equity Holdings Comparison sheet, shifting Mutual Fund/Status/
# Jan_2026/Dec_2025/Oct_2025/MoM/QoQ one column to the right (D:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C; Excel shifts C:I -> D:J and
# carries the header style/formatting along automatically.
$ws.Columns("C:C").Insert()

# Populate the new "Industry" column.
$industry = @{
    1  = "Industry"
    2  = "Pharmaceuticals & Biotechnology"
    3  = "Auto Components"
    4  = "Metals & Minerals Trading"
    5  = "Food Products"
    6  = "Industrial Products"
    7  = "Automobiles"
    8  = "Power"
    9  = "Electrical Equipment"
    10 = "Pharmaceuticals & Biotechnology"
    11 = "Finance"
    12 = "Power"
    13 = "Cement & Cement Products"
    14 = "Textiles & Apparels"
    15 = "Textiles & Apparels"
    16 = "Chemicals & Petrochemicals"
    17 = "Fertilizers & Agrochemicals"
    18 = "Pharmaceuticals & Biotechnology"
    19 = "Telecom - Services"
    20 = "IT - Software"
    21 = "Aerospace & Defense"
    22 = "Chemicals & Petrochemicals"
    23 = "Retailing"
    24 = "Pharmaceuticals & Biotechnology"
    25 = "Electrical Equipment"
    26 = "Chemicals & Petrochemicals"
}

foreach ($row in $industry.Keys) {
    $ws.Cells.Item($row, 3).Value = $industry[$row]
}
